$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous data range before writing the refreshed dataset
$ws.Range("A2:C49").ClearContents()

$rows = @(
    @('Tesla Model X I · Рестайлинг', 'https://cars.av.by/tesla/model-x/110458856', 282147),
    @('Tesla Model X I · Рестайлинг', 'https://cars.av.by/tesla/model-x/109287918', 270985),
    @('Tesla Model X I · Рестайлинг', 'https://cars.av.by/tesla/model-x/109256994', 270670),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/103694589', 256642),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/104400486', 253135),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/107287543', 251860),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/108444907', 223164),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/105562879', 216788),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109631463', 207227),
    @('Tesla Model X I · Рестайлинг, 6 мест', 'https://cars.av.by/tesla/model-x/109872261', 207227),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/109578074', 192880),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110405384', 191286),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/110121889', 191286),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/109905270', 191283),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110405848', 191254),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/106217374', 190967),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/100804369', 188098),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109599275', 187141),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110405800', 159373),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/101410920', 159086),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/109699442', 156217),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/110141975', 156214),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109158466', 146621),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/107479694', 146334),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/106159673', 143465),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/104415641', 143465),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/103933677', 143465),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/103900859', 141870),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/110157068', 137088),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/107218202', 135494),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/108738548', 133581),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/108904965', 132306),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/107436493', 131031),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/104192514', 127524),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/109356421', 127492),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/104745779', 125930),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109152501', 122742),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/110353371', 122742),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/106453955', 120829),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/107295316', 119554),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/107376045', 119554),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/110399301', 117957),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/108229841', 114453),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/106091296', 111902),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109805552', 111584),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/108427822', 111584),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/106624429', 111265),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/106319825', 111265),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110406206', 108392),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/108365228', 103613),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/109126708', 102016),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/109303397', 101700),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/108439391', 100425),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110225839', 98831),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/110297676', 98672),
    @('Tesla Model X I', 'https://cars.av.by/tesla/model-x/110432057', 98512),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/104088131', 95640),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/108658646', 95324),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/107709262', 92455),
    @('Tesla Model X I, 7 мест', 'https://cars.av.by/tesla/model-x/109699771', 90861),
    @('Tesla Model X I, 6 мест', 'https://cars.av.by/tesla/model-x/110441062', 87673)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

"Done writing $($rows.Count) rows. UsedRange: " + $ws.UsedRange.Address()
